$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29: MSRP bump
$ws.Range("D29").Value = 53100

# Row 30: MSRP bump
$ws.Range("D30").Value = 55890

# Row 31: MSRP bump
$ws.Range("D31").Value = 64365

# Row 32: MSRP bump + DPHF correction
$ws.Range("D32").Value = 86580
$ws.Range("E32").Value = 1025

# Row 33: MSRP bump + DPHF correction
$ws.Range("D33").Value = 91580
$ws.Range("E33").Value = 1025

# Row 34: fill in previously blank MSRP with real number (matches format of D29:D33)
$ws.Range("D34").NumberFormat = $ws.Range("D33").NumberFormat
$ws.Range("D34").Value = 99310
$ws.Range("E34").Value = 1025

# Match the saved selection
$ws.Range("D29").Select()
